$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.349205374717712
$ws.Range("B1").Value = 2.383350372314453
$ws.Range("C1").Value = 2.997998476028442
$ws.Range("D1").Value = 3.449360847473145
$ws.Range("E1").Value = 1.192476153373718
